$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (2023-10-03) for every
# data row from row 2 through row 140. Bump each of those cells by one day
# to 2023-10-04, matching the diff which changes 45202 -> 45203.
foreach ($row in 2..140) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = $cell.Value2 + 1
}
